# Updated cryptos list on Sun Mar 26 21:27:30 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns keep their original text formatting
# so numeric-looking strings (e.g. "0.3540", "92.80", "0.00001060") are not
# auto-converted to numbers (which would strip trailing zeros / use sci notation).
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "27.844.65"
$ws.Range("E2").Value = "  +1.40%  "

# Row 3
$ws.Range("D3").Value = "1.762.10"
$ws.Range("E3").Value = "  +1.43%  "

# Row 4
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.16%  "

# Row 5
$ws.Range("D5").Value = "327.89"
$ws.Range("E5").Value = "  +1.58%  "

# Row 6
$ws.Range("E6").Value = "  +0.11%  "

# Row 7
$ws.Range("D7").Value = "0.4462"
$ws.Range("E7").Value = "  -2.09%  "

# Row 8
$ws.Range("D8").Value = "0.3540"
$ws.Range("E8").Value = "  +0.36%  "

# Row 9
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.07404"
$ws.Range("E9").Value = "  +0.19%  "

# Row 10
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").Value = "41.94"
$ws.Range("E10").Value = "  +1.65%  "

# Row 11
$ws.Range("D11").Value = "1.098"
$ws.Range("E11").Value = "  +2.26%  "

# Row 12
$ws.Range("E12").Value = "  +0.13%  "

# Row 13
$ws.Range("D13").Value = "20.87"
$ws.Range("E13").Value = "  +2.30%  "

# Row 14
$ws.Range("D14").Value = "6.016"
$ws.Range("E14").Value = "  +1.79%  "

# Row 15
$ws.Range("D15").Value = "7.231"
$ws.Range("E15").Value = "  +2.69%  "

# Row 16
$ws.Range("D16").Value = "1.763.88"
$ws.Range("E16").Value = "  +1.40%  "

# Row 17
$ws.Range("D17").Value = "92.80"
$ws.Range("E17").Value = "  +1.81%  "

# Row 18
$ws.Range("D18").Value = "0.00001060"
$ws.Range("E18").Value = "  +0.63%  "

# Row 19
$ws.Range("D19").Value = "0.06428"
$ws.Range("E19").Value = "  +1.37%  "

# Row 21
$ws.Range("D21").Value = "17.07"
$ws.Range("E21").Value = "  +2.92%  "

# Row 22
$ws.Range("D22").Value = "5.753"
$ws.Range("E22").Value = "  +0.58%  "

# Row 23
$ws.Range("D23").Value = "27.893.04"
$ws.Range("E23").Value = "  +1.46%  "

# Row 24
$ws.Range("D24").Value = "11.23"
$ws.Range("E24").Value = "  +1.03%  "

# Row 25
$ws.Range("D25").Value = "2.109"
$ws.Range("E25").Value = "  +1.69%  "

# Row 26
$ws.Range("D26").Value = "160.85"
$ws.Range("E26").Value = "  -0.52%  "

# Row 27
$ws.Range("D27").Value = "20.36"
$ws.Range("E27").Value = "  +1.92%  "

# Row 28
$ws.Range("D28").Value = "1.967.48"
$ws.Range("E28").Value = "  +1.80%  "

# Row 29
$ws.Range("D29").Value = "2.142"
$ws.Range("E29").Value = "  +4.79%  "

# Row 30
$ws.Range("D30").Value = "124.25"
$ws.Range("E30").Value = "  -0.37%  "

# Row 31
$ws.Range("D31").Value = "1.098"
$ws.Range("E31").Value = "  +5.28%  "

# Row 32
$ws.Range("D32").Value = "0.09194"
$ws.Range("E32").Value = "  +1.22%  "

# Row 33
$ws.Range("D33").Value = "5.643"
$ws.Range("E33").Value = "  +4.89%  "

# Row 34
$ws.Range("E34").Value = "  +1.20%  "

# Row 35
$ws.Range("D35").Value = "11.83"
$ws.Range("E35").Value = "  +2.06%  "

# Row 36
$ws.Range("D36").Value = "0.06183"
$ws.Range("E36").Value = "  +3.90%  "

# Row 37
$ws.Range("D37").Value = "0.02278"
$ws.Range("E37").Value = "  +0.51%  "

# Row 38
$ws.Range("D38").Value = "0.2098"
$ws.Range("E38").Value = "  +2.03%  "

# Row 39
$ws.Range("D39").Value = "0.6303"
$ws.Range("E39").Value = "  +1.34%  "

# Row 40
$ws.Range("D40").Value = "4.942"
$ws.Range("E40").Value = "  +1.37%  "

# Row 41
$ws.Range("D41").Value = "1.181"
$ws.Range("E41").Value = "  -0.80%  "

# Row 42
$ws.Range("D42").Value = "1.395"
$ws.Range("E42").Value = "  +1.62%  "

# Row 43
$ws.Range("D43").Value = "7.859"
$ws.Range("E43").Value = "  +2.30%  "

# Row 44
$ws.Range("D44").Value = "13.29"
$ws.Range("E44").Value = "  +1.51%  "

# Row 45
$ws.Range("D45").Value = "3.742"
$ws.Range("E45").Value = "  +1.12%  "

# Row 46
$ws.Range("D46").Value = "0.5848"
$ws.Range("E46").Value = "  +1.06%  "

# Row 47
$ws.Range("D47").Value = "122.18"
$ws.Range("E47").Value = "  +0.33%  "

# Row 48
$ws.Range("D48").Value = "1.950"
$ws.Range("E48").Value = "  +1.31%  "

# Row 49
$ws.Range("D49").Value = "0.06888"
$ws.Range("E49").Value = "  +0.74%  "

# Row 50
$ws.Range("D50").Value = "1.133"
$ws.Range("E50").Value = "  +1.84%  "

# Row 51
$ws.Range("D51").Value = "72.74"
$ws.Range("E51").Value = "  +2.39%  "
